# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> the "Integral" theme, used by the (single) Slide Master
#   ppt/theme/theme2.xml  -> the default "Office Theme", used by the Notes Master
#
# The authored edit swaps the two themes' colour palettes: the Slide Master
# (theme1.xml) is switched over to the stock "Office Theme" colour scheme
# (the palette that used to live only in theme2.xml / the Notes Master).
#
# PowerPoint's object model exposes theme colour edits through
# ThemeColorScheme.Colors(i).RGB (indices 1-12, in Dark1/Light1/Dark2/Light2/
# Accent1-6/Hyperlink/FollowedHyperlink order) -- that's the supported,
# scriptable surface for "apply a different colour theme" (the font scheme
# and format scheme are already identical between the two theme parts in
# this deck, so only the colours need to change).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Target palette: the stock Office Theme colours (was previously only on
# theme2.xml / the Notes Master).
$colors.Colors(1).RGB  = 0         # Dark 1    - 000000
$colors.Colors(2).RGB  = 16777215  # Light 1   - FFFFFF
$colors.Colors(3).RGB  = 6968388   # Dark 2    - 44546A
$colors.Colors(4).RGB  = 15132391  # Light 2   - E7E6E6
$colors.Colors(5).RGB  = 13998939  # Accent 1  - 5B9BD5
$colors.Colors(6).RGB  = 3243501   # Accent 2  - ED7D31
$colors.Colors(7).RGB  = 10855845  # Accent 3  - A5A5A5
$colors.Colors(8).RGB  = 49407     # Accent 4  - FFC000
$colors.Colors(9).RGB  = 12874308  # Accent 5  - 4472C4
$colors.Colors(10).RGB = 4697456   # Accent 6  - 70AD47
$colors.Colors(11).RGB = 12673797  # Hyperlink - 0563C1
$colors.Colors(12).RGB = 7491477   # Followed Hyperlink - 954F72
